$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to reflect the random-state comparison being added.
$ws.Name = "randomstate = 1"

# Clear the now-unused "divider" styles (no real border, just font/border
# flags left over from earlier formatting passes) off the raw-metrics table.
$ws.Range("A15:A32").Style = "Normal"
$ws.Range("B12").Style = "Normal"

# The trailing blank, styled row is no longer needed - remove it so the
# used range shrinks back down to the real data (A1:G32).
$ws.Rows(33).Delete()

# Widen column A slightly to fit the (now finalized) row labels.
$ws.Columns("A").ColumnWidth = 11.63

# Leave the selection where the author was last working.
$ws.Range("B26:F26").Select()
